# Removendo o campo DATA dos acervos audiovisuais
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G contains the "Data" header/field; delete the whole column so
# everything to its right shifts left by one.
$ws.Range("G1").EntireColumn.Delete()

$ws.Range("I8").Select()
